# The underlying records for rows 9-13 (data rows, row 1 is the header)
# have been re-keyed: each row's entire contents now match a different
# original row, per this permutation (new row -> source old row):
#   9 -> 10, 10 -> 13, 11 -> 12, 12 -> 9, 13 -> 11
#
# Implement this as a full-row value swap using the Excel object model:
# snapshot every row's values first (so overwrites don't clobber a value
# still needed by a later assignment), then write them back according to
# the permutation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 9
$lastRow = 13
$firstCol = "A"
$lastCol = "AY"

# Columns Y and AA store dates as plain text (e.g. "2023-08-11"); force
# text format up front so re-assigning those values through .Value does
# not get reinterpreted as a date serial number by Excel.
$ws.Range("Y$firstRow`:Y$lastRow").NumberFormat = "@"
$ws.Range("AA$firstRow`:AA$lastRow").NumberFormat = "@"

# Snapshot current values of each row in the affected range.
$rowValues = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowValues[$r] = $ws.Range("$firstCol$r`:$lastCol$r").Value()
}

# new row r gets the content that used to live in row $sourceRow[r]
$sourceRow = @{
    9  = 10
    10 = 13
    11 = 12
    12 = 9
    13 = 11
}

foreach ($r in 9..13) {
    $src = $sourceRow[$r]
    $ws.Range("$firstCol$r`:$lastCol$r").Value = $rowValues[$src]
}
